$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2316021.2
$ws.Range("J17").Value = 2316021.2
$ws.Range("L17").Value = 6948063.600000001
$ws.Range("N17").Value = -6948399.600000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1839.8
$ws.Range("I19").Value = 2658.2
$ws.Range("K19").Value = 2658.2
$ws.Range("M19").Value = -2483.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 298.81482
$ws.Range("I33").Value = 314.56
$ws.Range("J33").Value = 102
$ws.Range("K33").Value = 314.56
$ws.Range("L33").Value = 102
$ws.Range("M33").Value = -85.56
$ws.Range("N33").Value = -560

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 28475.5
$ws.Range("I43").Value = 6500
$ws.Range("J43").Value = 50451
$ws.Range("K43").Value = 6500
$ws.Range("L43").Value = 50451
$ws.Range("M43").Value = -6431
$ws.Range("N43").Value = -50589

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 1000
$ws.Range("J45").Value = 1000
$ws.Range("L45").Value = 3000
$ws.Range("N45").Value = -3384

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 32390
$ws.Range("J123").Value = 32390
$ws.Range("L123").Value = 32390
$ws.Range("N123").Value = -42190

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 32620
$ws.Range("J140").Value = 32620
$ws.Range("L140").Value = 32620
$ws.Range("N140").Value = -42980

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9805.893
$ws.Range("I32").Value = 8973.263999999999
$ws.Range("K32").Value = 8973.263999999999
$ws.Range("M32").Value = -8686.263999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 17858726
$ws.Range("I61").Value = 20834846
$ws.Range("K61").Value = 20834846
$ws.Range("M61").Value = -20834634

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 55557644
$ws.Range("I63").Value = 62501852
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 62501852
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -62501166
$ws.Range("N63").Value = -5372

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 55557644
$ws.Range("I66").Value = 62501852
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 312509260
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -312505828
$ws.Range("N66").Value = -26864

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 25003120
$ws.Range("I74").Value = 41667988
$ws.Range("J74").Value = 5817.75
$ws.Range("K74").Value = 41667988
$ws.Range("L74").Value = 5817.75
$ws.Range("M74").Value = -41667114
$ws.Range("N74").Value = -7565.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 25003120
$ws.Range("I77").Value = 41667988
$ws.Range("J77").Value = 5817.75
$ws.Range("K77").Value = 208339940
$ws.Range("L77").Value = 29088.75
$ws.Range("M77").Value = -208335572
$ws.Range("N77").Value = -37824.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1459.5
$ws.Range("I110").Value = 785.5
$ws.Range("J110").Value = 3481.5
$ws.Range("K110").Value = 785.5
$ws.Range("L110").Value = 3481.5
$ws.Range("M110").Value = 1259.5
$ws.Range("N110").Value = -7571.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7814257.5
$ws.Range("I132").Value = 10418085
$ws.Range("J132").Value = 2774.75
$ws.Range("K132").Value = 31254255
$ws.Range("L132").Value = 8324.25
$ws.Range("M132").Value = -31251725
$ws.Range("N132").Value = -13384.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 17858726
$ws.Range("I136").Value = 20834846
$ws.Range("K136").Value = 62504538
$ws.Range("M136").Value = -62501988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2688.75
$ws.Range("I20").Value = 3630
$ws.Range("J20").Value = 1120
$ws.Range("K20").Value = 3630
$ws.Range("L20").Value = 1120
$ws.Range("M20").Value = -3383
$ws.Range("N20").Value = -1614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 16085.5
$ws.Range("I82").Value = 15128.5
$ws.Range("J82").Value = 17042.5
$ws.Range("K82").Value = 15128.5
$ws.Range("L82").Value = 17042.5
$ws.Range("M82").Value = -14745.5
$ws.Range("N82").Value = -17808.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 16085.5
$ws.Range("I85").Value = 15128.5
$ws.Range("J85").Value = 17042.5
$ws.Range("K85").Value = 15128.5
$ws.Range("L85").Value = 17042.5
$ws.Range("M85").Value = -13802.5
$ws.Range("N85").Value = -19694.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 820.8823
$ws.Range("I94").Value = 720
$ws.Range("J94").Value = 965
$ws.Range("K94").Value = 720
$ws.Range("L94").Value = 965
$ws.Range("M94").Value = -269
$ws.Range("N94").Value = -1867

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11116401
$ws.Range("I31").Value = 5617.5557
$ws.Range("J31").Value = 111113450
$ws.Range("K31").Value = 5617.5557
$ws.Range("L31").Value = 111113450
$ws.Range("M31").Value = -5322.5557
$ws.Range("N31").Value = -111114040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 11116401
$ws.Range("I34").Value = 5617.5557
$ws.Range("J34").Value = 111113450
$ws.Range("K34").Value = 5617.5557
$ws.Range("L34").Value = 111113450
$ws.Range("M34").Value = -5415.5557
$ws.Range("N34").Value = -111113854

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1383.6364
$ws.Range("I99").Value = 1317.8948
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 1317.8948
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = 180.1052
$ws.Range("N99").Value = -4796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1383.6364
$ws.Range("I126").Value = 1317.8948
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 3953.6844
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -1483.6844
$ws.Range("N126").Value = -10340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 596315.8
$ws.Range("I134").Value = 1068.5385
$ws.Range("J134").Value = 1701775.1
$ws.Range("K134").Value = 3205.6155
$ws.Range("L134").Value = 5105325.300000001
$ws.Range("M134").Value = -670.6155000000003
$ws.Range("N134").Value = -5110395.300000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 40757.5
$ws.Range("J140").Value = 40757.5
$ws.Range("L140").Value = 40757.5
$ws.Range("N140").Value = -51117.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2470.6
$ws.Range("J80").Value = 2400.75
$ws.Range("L80").Value = 7202.25
$ws.Range("N80").Value = -9074.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 2470.6
$ws.Range("J83").Value = 2400.75
$ws.Range("L83").Value = 21606.75
$ws.Range("N83").Value = -30966.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 2737.5
$ws.Range("J105").Value = 2737.5
$ws.Range("L105").Value = 8212.5
$ws.Range("N105").Value = -13454.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10104432
$ws.Range("I80").Value = 13891581
$ws.Range("J80").Value = 5367.222
$ws.Range("K80").Value = 13891581
$ws.Range("L80").Value = 5367.222
$ws.Range("M80").Value = -13890583
$ws.Range("N80").Value = -7363.222

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 10104432
$ws.Range("I83").Value = 13891581
$ws.Range("J83").Value = 5367.222
$ws.Range("K83").Value = 69457905
$ws.Range("L83").Value = 26836.11
$ws.Range("M83").Value = -69452913
$ws.Range("N83").Value = -36820.11

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 23656.445
$ws.Range("J123").Value = 23656.445
$ws.Range("L123").Value = 23656.445
$ws.Range("N123").Value = -28556.445

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3347.457
$ws.Range("I132").Value = 2348.3635
$ws.Range("K132").Value = 7045.0905
$ws.Range("M132").Value = -4515.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 49547
$ws.Range("J138").Value = 49547
$ws.Range("L138").Value = 49547
$ws.Range("N138").Value = -59827

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1899.8572
$ws.Range("I61").Value = 1833.1666
$ws.Range("K61").Value = 1833.1666
$ws.Range("M61").Value = -1631.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1899.8572
$ws.Range("I113").Value = 1833.1666
$ws.Range("K113").Value = 1833.1666
$ws.Range("M113").Value = 336.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 15160879
$ws.Range("I132").Value = 7149.857
$ws.Range("J132").Value = 26326786
$ws.Range("K132").Value = 21449.571
$ws.Range("L132").Value = 78980358
$ws.Range("M132").Value = -18919.571
$ws.Range("N132").Value = -78985418

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 60663.57
$ws.Range("J139").Value = 60663.57
$ws.Range("L139").Value = 60663.57
$ws.Range("N139").Value = -70943.57000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1905.125
$ws.Range("I81").Value = 279.66666
$ws.Range("J81").Value = 2880.4
$ws.Range("K81").Value = 559.33332
$ws.Range("L81").Value = 5760.8
$ws.Range("M81").Value = 501.66668
$ws.Range("N81").Value = -7882.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1905.125
$ws.Range("I84").Value = 279.66666
$ws.Range("J84").Value = 2880.4
$ws.Range("K84").Value = 2796.6666
$ws.Range("L84").Value = 28804
$ws.Range("M84").Value = 2507.3334
$ws.Range("N84").Value = -39412

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1468.9286
$ws.Range("I132").Value = 689
$ws.Range("K132").Value = 2067
$ws.Range("M132").Value = 463
